$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mojibake (UTF-8 double-encoding) in the Regional Economic Communities note (A103)
$ws.Range("A103").Value = 'Regional Economic Communities: CEN-SAD = "Community of Sahel-Saharan States"; COMESA = "Common Market for Eastern and Southern Africa"; EAC = "East African Community"; ECCAS = "Economic Community of Central African States"; ECOWAS = "Economic Community of West African States"; IGAD = "Intergovernmental Authority on Development"; SADC = "Southern African Development Community"; UMA = "Arab Maghreb Union"; PALOP = "Países Africanos de Língua Oficial Portuguesa"; ASEAN = "Association of Southeast Asian Nations"; MERCOSUR = "Mercado Común del Sur". EU27 = "European Union (27 members)". OECD = "Organisation for Economic Co-operation and Development".'

# Updated data values (columns G, H, M, N) for several country/aggregate rows
$ws.Range("G13").Value = 3.6880769999999998
$ws.Range("H13").Value = 5.9316180000000003
$ws.Range("M13").Value = 2.604393
$ws.Range("N13").Value = 4.9021100000000004
$ws.Range("G23").Value = 3.1000480000000001
$ws.Range("H23").Value = 6.7829569999999997
$ws.Range("M23").Value = 5.8668360000000002
$ws.Range("N23").Value = 13.164300000000001
$ws.Range("G61").Value = 3.8492820000000001
$ws.Range("H61").Value = 7.7195770000000001
$ws.Range("M61").Value = 7.7740309999999999
$ws.Range("N61").Value = 15.644133
$ws.Range("G62").Value = 5.0606439999999999
$ws.Range("H62").Value = 11.552286
$ws.Range("M62").Value = 3.8447619999999998
$ws.Range("N62").Value = 8.3914790000000004
$ws.Range("G63").Value = 5.0128539999999999
$ws.Range("H63").Value = 10.070586
$ws.Range("M63").Value = 1.3720220000000001
$ws.Range("N63").Value = 2.7035079999999998
$ws.Range("G64").Value = 2.8322080000000001
$ws.Range("H64").Value = 7.5359119999999997
$ws.Range("M64").Value = 0.76941300000000001
$ws.Range("N64").Value = 2.0728469999999999
$ws.Range("G65").Value = 8.6177609999999998
$ws.Range("H65").Value = 16.810333
$ws.Range("M65").Value = 3.4092060000000002
$ws.Range("N65").Value = 6.7084599999999996
$ws.Range("G66").Value = 5.0155060000000002
$ws.Range("H66").Value = 10.15279
$ws.Range("M66").Value = 1.5092080000000001
$ws.Range("N66").Value = 3.0190730000000001
$ws.Range("G67").Value = 6.3223779999999996
$ws.Range("H67").Value = 15.892503
$ws.Range("M67").Value = 2.323369
$ws.Range("N67").Value = 6.4273429999999996
$ws.Range("G68").Value = 6.2897109999999996
$ws.Range("H68").Value = 14.574522
$ws.Range("M68").Value = 4.1936600000000004
$ws.Range("N68").Value = 8.5402509999999996
$ws.Range("G70").Value = 2.6094020000000002
$ws.Range("H70").Value = 6.0590979999999997
$ws.Range("M70").Value = 4.8048330000000004
$ws.Range("N70").Value = 11.402274999999999
$ws.Range("G71").Value = 3.8492820000000001
$ws.Range("H71").Value = 7.7195770000000001
$ws.Range("M71").Value = 7.7740309999999999
$ws.Range("N71").Value = 15.644133
$ws.Range("G73").Value = 3.265533
$ws.Range("H73").Value = 5.4695780000000003
$ws.Range("M73").Value = 3.5521449999999999
$ws.Range("N73").Value = 7.6465860000000001
$ws.Range("G75").Value = 1.9722109999999999
$ws.Range("H75").Value = 4.9169499999999999
$ws.Range("M75").Value = 3.9713539999999998
$ws.Range("N75").Value = 9.4490379999999998
$ws.Range("N76").Value = 6.1687029999999998
$ws.Range("G77").Value = 2.7015150000000001
$ws.Range("H77").Value = 7.2976869999999998
$ws.Range("M77").Value = 0.58216599999999996
$ws.Range("N77").Value = 1.6254580000000001
$ws.Range("G78").Value = 3.4954329999999998
$ws.Range("H78").Value = 5.3012680000000003
$ws.Range("M78").Value = 0.156859
$ws.Range("N78").Value = 0.174072
$ws.Range("G79").Value = 2.4443950000000001
$ws.Range("H79").Value = 4.9907579999999996
$ws.Range("M79").Value = 0.076133000000000006
$ws.Range("N79").Value = 0.16664599999999999
$ws.Range("G80").Value = 4.3597000000000001
$ws.Range("H80").Value = 12.657641999999999
$ws.Range("M80").Value = 1.413343
$ws.Range("N80").Value = 3.2451479999999999
$ws.Range("G81").Value = 5.067717
$ws.Range("H81").Value = 13.947978000000001
$ws.Range("M81").Value = 0.24724699999999999
$ws.Range("N81").Value = 0.43498500000000001
$ws.Range("G82").Value = 5.1846740000000002
$ws.Range("H82").Value = 11.356698
$ws.Range("M82").Value = 4.2749940000000004
$ws.Range("N82").Value = 9.3021049999999992
$ws.Range("G83").Value = 5.0076239999999999
$ws.Range("H83").Value = 9.7009749999999997
$ws.Range("M83").Value = 1.479241
$ws.Range("N83").Value = 2.9197540000000002
$ws.Range("G84").Value = 1.910407
$ws.Range("H84").Value = 5.5858299999999996
$ws.Range("M84").Value = 6.9542669999999998
$ws.Range("N84").Value = 17.518986000000002
$ws.Range("G86").Value = 5.8027290000000002
$ws.Range("H86").Value = 13.070997
$ws.Range("M86").Value = 4.3318209999999997
$ws.Range("N86").Value = 8.7714800000000004
$ws.Range("G88").Value = 5.0078449999999997
$ws.Range("H88").Value = 11.218268999999999
$ws.Range("M88").Value = 0.52309499999999998
$ws.Range("N88").Value = 0.81241300000000005
$ws.Range("G89").Value = 7.7916840000000001
$ws.Range("H89").Value = 15.276524999999999
$ws.Range("M89").Value = 1.939792
$ws.Range("N89").Value = 3.7680899999999999
$ws.Range("G90").Value = 2.66717
$ws.Range("H90").Value = 5.5722719999999999
$ws.Range("M90").Value = 0.066990999999999995
$ws.Range("N90").Value = 0.083710999999999994
$ws.Range("G91").Value = 1.819202
$ws.Range("H91").Value = 4.9321520000000003
$ws.Range("M91").Value = 6.5052770000000004
$ws.Range("N91").Value = 16.651022000000001
$ws.Range("N92").Value = 17.987169000000002
$ws.Range("G93").Value = 4.7179250000000001
$ws.Range("H93").Value = 8.0905889999999996
$ws.Range("M93").Value = 1.6488719999999999
$ws.Range("N93").Value = 3.7966880000000001
$ws.Range("G94").Value = 2.6661510000000002
$ws.Range("H94").Value = 10.022695000000001
$ws.Range("M94").Value = 1.3225880000000001
$ws.Range("N94").Value = 2.4069240000000001
$ws.Range("G95").Value = 1.5134860000000001
$ws.Range("H95").Value = 3.3022170000000002
$ws.Range("M95").Value = 7.4495909999999999
$ws.Range("N95").Value = 18.180109999999999
$ws.Range("G96").Value = 6.4817330000000002
$ws.Range("H96").Value = 12.34811
$ws.Range("M96").Value = 2.3831570000000002
$ws.Range("N96").Value = 5.2066929999999996
$ws.Range("G97").Value = 2.9745710000000001
$ws.Range("H97").Value = 6.6470349999999998
$ws.Range("M97").Value = 6.6437109999999997
$ws.Range("N97").Value = 14.665988
$ws.Range("G98").Value = 5.4976770000000004
$ws.Range("H98").Value = 13.739038000000001
$ws.Range("M98").Value = 4.3044279999999997
$ws.Range("N98").Value = 9.1160019999999999
